$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("azúcar", 3753.876470588235),
    @("huevos", 2679.706666666666),
    @("leche semidesnatada", 21982.638),
    @("aceite vegetal", 10260.66304347826),
    @("cacao en polvo 32%", "No prices found"),
    @("harina 1kg", 1474.8125),
    @("polvo de hornear", 5656.67),
    @("sal", 3033.001034482759),
    @("margarina", 7329.649230769231),
    @("leche condensada", 4797.56875),
    @("crema de leche", 6428.746551724138),
    @("cacao en polvo 50%", "No prices found"),
    @("granulado de chocolate", 9761.155999999999)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
